# Auto-generated edit script: updates FFXIV leve market-price snapshot values
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW) to match the
# latest scheduled-runner price pull. Cells that no longer have a meaningful
# profit figure are cleared; cells that gain one are written for the first time.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 513.9
$ws.Range("I15").Value = 513.9
$ws.Range("K15").Value = 1541.7
$ws.Range("M15").Value = -1372.7
# Row 18
$ws.Range("H18").Value = 2113.2
$ws.Range("J18").Value = 3783
$ws.Range("L18").Value = 3783
$ws.Range("N18").Value = -4351
# Row 19
$ws.Range("H19").Value = 1785.1428
$ws.Range("I19").Value = 1749.25
$ws.Range("K19").Value = 1749.25
$ws.Range("M19").Value = -1574.25
# Row 38
$ws.Range("H38").Value = 18.25
$ws.Range("I38").Value = 18.25
$ws.Range("K38").Value = 54.75
$ws.Range("M38").Value = 317.25
# Row 42
$ws.Range("H42").Value = 131.5
$ws.Range("J42").Value = 175
$ws.Range("L42").Value = 525
$ws.Range("N42").Value = -985
# Row 55
$ws.Range("H55").Value = 425.57144
$ws.Range("I55").Value = 328.66666
$ws.Range("J55").Value = 600
$ws.Range("K55").Value = 328.66666
$ws.Range("L55").Value = 600
$ws.Range("M55").Value = -114.66666
$ws.Range("N55").Value = -1028
# Row 62
$ws.Range("H62").Value = 7249.25
$ws.Range("J62").Value = 12500
$ws.Range("L62").Value = 12500
$ws.Range("N62").Value = -13748
# Row 65
$ws.Range("H65").Value = 7249.25
$ws.Range("J65").Value = 12500
$ws.Range("L65").Value = 62500
$ws.Range("N65").Value = -68740

$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 750
$ws.Range("I88").Value = 750
$ws.Range("K88").Value = 750
$ws.Range("M88").Value = -344
# Row 91
$ws.Range("H91").Value = 750
$ws.Range("I91").Value = 750
$ws.Range("K91").Value = 750
$ws.Range("M91").Value = 654
# Row 102
$ws.Range("H102").Value = 680
$ws.Range("I102").Value = 680
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 680
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 942
$ws.Range("N102").ClearContents()
# Row 132
$ws.Range("H132").Value = 3111.7
$ws.Range("I132").Value = 1236.1666
$ws.Range("J132").Value = 5925
$ws.Range("K132").Value = 3708.4998
$ws.Range("L132").Value = 17775
$ws.Range("M132").Value = -1178.4998
$ws.Range("N132").Value = -22835

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4759.15
$ws.Range("I86").Value = 4815.4614
$ws.Range("J86").Value = 4654.5713
$ws.Range("K86").Value = 4815.4614
$ws.Range("L86").Value = 4654.5713
$ws.Range("M86").Value = -3692.4614
$ws.Range("N86").Value = -6900.5713
# Row 89
$ws.Range("H89").Value = 4759.15
$ws.Range("I89").Value = 4815.4614
$ws.Range("J89").Value = 4654.5713
$ws.Range("K89").Value = 24077.307
$ws.Range("L89").Value = 23272.8565
$ws.Range("M89").Value = -18461.307
$ws.Range("N89").Value = -34504.85649999999

$ws = $wb.Worksheets.Item("CRP")
# Row 60
$ws.Range("H60").Value = 30051
$ws.Range("J60").Value = 30103
$ws.Range("L60").Value = 30103
$ws.Range("N60").Value = -31125
# Row 69
$ws.Range("H69").Value = 38000
$ws.Range("I69").Value = 7000
$ws.Range("K69").Value = 7000
$ws.Range("M69").Value = -6251
# Row 72
$ws.Range("H72").Value = 38000
$ws.Range("I72").Value = 7000
$ws.Range("K72").Value = 21000
$ws.Range("M72").Value = -17256
# Row 95
$ws.Range("H95").Value = 13197
$ws.Range("J95").Value = 13197
$ws.Range("L95").Value = 13197
$ws.Range("N95").Value = -18689
# Row 103
$ws.Range("H103").Value = 33999.5
$ws.Range("I103").Value = 33999.5
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 33999.5
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -32827.5
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
# Row 12
$ws.Range("H12").Value = 209
$ws.Range("J12").Value = 335.75
$ws.Range("L12").Value = 1007.25
$ws.Range("N12").Value = -1353.25
# Row 112
$ws.Range("H112").Value = 2000
$ws.Range("J112").Value = 2000
$ws.Range("L112").Value = 6000
$ws.Range("N112").Value = -8216
# Row 131
$ws.Range("H131").Value = 1491.25
$ws.Range("J131").Value = 1733
$ws.Range("L131").Value = 5199
$ws.Range("N131").Value = -15279
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
# Row 6
$ws.Range("H6").Value = 4083.3333
$ws.Range("I6").Value = 1125
$ws.Range("J6").Value = 10000
$ws.Range("K6").Value = 1125
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = -1012
$ws.Range("N6").Value = -10226
# Row 7
$ws.Range("H7").Value = 2505000
$ws.Range("J7").Value = 5000000
$ws.Range("L7").Value = 5000000
$ws.Range("N7").Value = -5000224
# Row 8
$ws.Range("H8").Value = 2505000
$ws.Range("J8").Value = 5000000
$ws.Range("L8").Value = 5000000
$ws.Range("N8").Value = -5000278
# Row 10
$ws.Range("H10").Value = 75000
$ws.Range("I10").Value = 130000
$ws.Range("J10").Value = 20000
$ws.Range("K10").Value = 130000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = -129831
$ws.Range("N10").Value = -20338
# Row 11
$ws.Range("H11").Value = 14467000
$ws.Range("I11").Value = 15154231
$ws.Range("J11").Value = 10000000
$ws.Range("K11").Value = 15154231
$ws.Range("L11").Value = 10000000
$ws.Range("M11").Value = -15154092
$ws.Range("N11").Value = -10000278
# Row 13
$ws.Range("H13").Value = 200
$ws.Range("I13").Value = 200
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 200
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -61
$ws.Range("N13").ClearContents()
# Row 14
$ws.Range("H14").Value = 180971.14
$ws.Range("J14").Value = 87266.336
$ws.Range("L14").Value = 87266.336
$ws.Range("N14").Value = -87602.336
# Row 16
$ws.Range("H16").Value = 4083.3333
$ws.Range("I16").Value = 1125
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 1125
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -875
$ws.Range("N16").Value = -10500
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
# Row 18
$ws.Range("H18").Value = 10000000
$ws.Range("I18").Value = 10000000
$ws.Range("K18").Value = 10000000
$ws.Range("M18").Value = -9999707
# Row 19
$ws.Range("H19").Value = 3997.6667
$ws.Range("I19").Value = 996.5
$ws.Range("J19").Value = 10000
$ws.Range("K19").Value = 996.5
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = -708.5
$ws.Range("N19").Value = -10576
# Row 21
$ws.Range("H21").Value = 525
$ws.Range("J21").Value = 525
$ws.Range("L21").Value = 525
$ws.Range("N21").Value = -871
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
# Row 27
$ws.Range("H27").Value = 90
$ws.Range("I27").Value = 90
$ws.Range("K27").Value = 90
$ws.Range("M27").Value = 76
# Row 30
$ws.Range("H30").Value = 525
$ws.Range("J30").Value = 525
$ws.Range("L30").Value = 525
$ws.Range("N30").Value = -735
# Row 33
$ws.Range("H33").Value = 5017
$ws.Range("I33").Value = 5017
$ws.Range("K33").Value = 5017
$ws.Range("M33").Value = -4765
# Row 52
$ws.Range("H52").Value = 4899
$ws.Range("J52").Value = 4899
$ws.Range("L52").Value = 4899
$ws.Range("N52").Value = -5417

$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 4091
$ws.Range("I9").Value = 132
$ws.Range("K9").Value = 132
$ws.Range("M9").Value = 92
# Row 17
$ws.Range("H17").Value = 3498.3333
$ws.Range("I17").Value = 2495
$ws.Range("K17").Value = 2495
$ws.Range("M17").Value = -2325
# Row 22
$ws.Range("H22").Value = 1028.4286
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 1199.8
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 1199.8
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -1789.8
# Row 27
$ws.Range("H27").Value = 1028.4286
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 1199.8
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 1199.8
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -1413.8
# Row 42
$ws.Range("H42").Value = 40000000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 40000000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 40000000
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -40001126
# Row 49
$ws.Range("H49").Value = 40000000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 40000000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 40000000
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -40000294
# Row 61
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4798
# Row 68
$ws.Range("H68").Value = 4000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 4000
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -5498
# Row 71
$ws.Range("H71").Value = 4000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 20000
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -27488
# Row 113
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830
